$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Value = "23-10-2025"
$ws.Range("B32").Value = "The price of gold in India today is ₹12,508 per gram for 24 karat gold, ₹11,465 per gram for 22 karat gold and ₹9,381 per gram for 18 karat gold (also called 999 gold)."
